$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.791298
$ws.Range("H2").Value = 14.373894
$ws.Range("I2").Value = 0.2539858212527056
$ws.Range("J2").Value = 0.2539858212527057
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.938300333333333
$ws.Range("N2").Value = 17.814901
$ws.Range("O2").Value = 0.1393304991144193
$ws.Range("P2").Value = 0.1393304991144193
$ws.Range("Q2").Value = 28.45216651049934
$ws.Range("R2").Value = 256.069498594494
$ws.Range("S2").Value = 0.03538797124312517
$ws.Range("T2").Value = 0.03538797124312518

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.791298
$ws.Range("H3").Value = 14.373894
$ws.Range("I3").Value = 0.2539858212527056
$ws.Range("J3").Value = 0.2539858212527057
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.82685766666667
$ws.Range("N3").Value = 32.480573
$ws.Range("O3").Value = 0.2540308502198431
$ws.Range("P3").Value = 0.2540308502198431
$ws.Range("Q3").Value = 51.87470148458467
$ws.Range("R3").Value = 466.872313361262
$ws.Range("S3").Value = 0.0645202341166099
$ws.Range("T3").Value = 0.06452023411660991

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.791298
$ws.Range("H4").Value = 14.373894
$ws.Range("I4").Value = 0.2539858212527056
$ws.Range("J4").Value = 0.2539858212527057
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 25.85508933333334
$ws.Range("N4").Value = 77.565268
$ws.Range("O4").Value = 0.6066386506657375
$ws.Range("P4").Value = 0.6066386506657376
$ws.Range("Q4").Value = 123.8794378126213
$ws.Range("R4").Value = 1114.914940313592
$ws.Range("S4").Value = 0.1540776158929706
$ws.Range("T4").Value = 0.1540776158929706

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 9.844169
$ws.Range("H5").Value = 29.532507
$ws.Range("I5").Value = 0.5218375788805928
$ws.Range("J5").Value = 0.5218375788805928
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.938300333333333
$ws.Range("N5").Value = 17.814901
$ws.Range("O5").Value = 0.1393304991144193
$ws.Range("P5").Value = 0.1393304991144193
$ws.Range("Q5").Value = 58.45763205408967
$ws.Range("R5").Value = 526.1186884868071
$ws.Range("S5").Value = 0.07270789032209315
$ws.Range("T5").Value = 0.07270789032209315

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 9.844169
$ws.Range("H6").Value = 29.532507
$ws.Range("I6").Value = 0.5218375788805928
$ws.Range("J6").Value = 0.5218375788805928
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 10.82685766666667
$ws.Range("N6").Value = 32.480573
$ws.Range("O6").Value = 0.2540308502198431
$ws.Range("P6").Value = 0.2540308502198431
$ws.Range("Q6").Value = 106.5814166096123
$ws.Range("R6").Value = 959.2327494865111
$ws.Range("S6").Value = 0.1325628438397014
$ws.Range("T6").Value = 0.1325628438397014

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 9.844169
$ws.Range("H7").Value = 29.532507
$ws.Range("I7").Value = 0.5218375788805928
$ws.Range("J7").Value = 0.5218375788805928
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 25.85508933333334
$ws.Range("N7").Value = 77.565268
$ws.Range("O7").Value = 0.6066386506657375
$ws.Range("P7").Value = 0.6066386506657376
$ws.Range("Q7").Value = 254.5218689074307
$ws.Range("R7").Value = 2290.696820166876
$ws.Range("S7").Value = 0.3165668447187982
$ws.Range("T7").Value = 0.3165668447187983

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.228964
$ws.Range("H8").Value = 12.686892
$ws.Range("I8").Value = 0.2241765998667015
$ws.Range("J8").Value = 0.2241765998667015
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 5.938300333333333
$ws.Range("N8").Value = 17.814901
$ws.Range("O8").Value = 0.1393304991144193
$ws.Range("P8").Value = 0.1393304991144193
$ws.Range("Q8").Value = 25.11285833085467
$ws.Range("R8").Value = 226.015724977692
$ws.Range("S8").Value = 0.03123463754920099
$ws.Range("T8").Value = 0.03123463754920099

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.228964
$ws.Range("H9").Value = 12.686892
$ws.Range("I9").Value = 0.2241765998667015
$ws.Range("J9").Value = 0.2241765998667015
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 10.82685766666667
$ws.Range("N9").Value = 32.480573
$ws.Range("O9").Value = 0.2540308502198431
$ws.Range("P9").Value = 0.2540308502198431
$ws.Range("Q9").Value = 45.78639130545734
$ws.Range("R9").Value = 412.077521749116
$ws.Range("S9").Value = 0.05694777226353174
$ws.Range("T9").Value = 0.05694777226353175

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.228964
$ws.Range("H10").Value = 12.686892
$ws.Range("I10").Value = 0.2241765998667015
$ws.Range("J10").Value = 0.2241765998667015
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 25.85508933333334
$ws.Range("N10").Value = 77.565268
$ws.Range("O10").Value = 0.6066386506657375
$ws.Range("P10").Value = 0.6066386506657376
$ws.Range("Q10").Value = 109.3402420074507
$ws.Range("R10").Value = 984.062178067056
$ws.Range("S10").Value = 0.1359941900539688
$ws.Range("T10").Value = 0.1359941900539688
